$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.946.89"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "2.679.90"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.106"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.15%  "
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.368"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.14%  "
$ws.Range("D13").Value = "3.154.87"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "62.808.82"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000148"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.91%  "
$ws.Range("D17").Value = "2.678.76"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.75%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.504"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").Value = "0.0₃0860"
$ws.Range("E28").Value = "  -6.13%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "341.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.954"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.81%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.50%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.46%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.21%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0565"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.617"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0973"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.15%  "
$ws.Range("D51").Value = "2.097.26"
$ws.Range("E51").Value = "  -1.70%  "
